# Implemented parallel bayesian with processor limits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits ---------------------------------------------------------

# Column F (rows 2-67): all values change from 5 to 0.5
$ws.Range("F2:F67").Value = 0.5

# Column E (rows 45-67): values change from 1E-3 (0.001) to 0.05
$ws.Range("E45:E67").Value = 0.05

# --- New column N: currency-formatted empty cell on row 23 -------------

$ws.Columns.Item(14).ColumnWidth = 14.28515625
$ws.Cells.Item(23, 14).Style = "Currency"

# --- Sheet view cosmetics (scroll position / selection) -----------------

$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("E8").Select()
